# This workbook's data rows were re-shuffled: the full content of a number of
# rows (2-93) was rotated among themselves (each destination row receives the
# previous full content of another row in the same permutation cycle). No new
# data values are introduced and no rows are added/removed - only existing row
# contents trade places. This script reproduces that reshuffle using the
# Excel COM object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..AY (1..51) hold the data for every record row.
$firstCol = "A"
$lastCol  = "AY"

# A handful of text columns contain values that *look* numeric or look like a
# date ("1", "2026-02-03", ...). Excel helpfully (but unhelpfully for us)
# auto-converts such strings to numbers/dates when they are assigned through
# the COM object model. Force those specific columns to stay plain text for
# the duration of the rewrite so the original string values round-trip
# unchanged.
$protectedCols = @("I", "Y", "AA")

# Each inner array below is one permutation cycle: the new content of
# cycle[i] is the old content of cycle[i+1] (wrapping around). Together the
# cycles cover every row whose data moved.
$cycles = @(
    ,@(2, 4)
    ,@(5, 6)
    ,@(9, 14, 11, 10, 13)
    ,@(15, 17, 16)
    ,@(18, 22, 20)
    ,@(19, 23, 21)
    ,@(31, 32)
    ,@(37, 38)
    ,@(47, 49, 50)
    ,@(54, 56, 58, 55)
    ,@(65, 66)
    ,@(77, 81)
    ,@(78, 82)
    ,@(83, 85, 86, 87, 88, 84)
    ,@(89, 91, 93)
    ,@(90, 92)
)

# Flatten to the full list of affected rows.
$allRows = @()
foreach ($cycle in $cycles) {
    foreach ($r in $cycle) {
        $allRows += $r
    }
}

# 1) Snapshot every affected row's full contents (A:AY) BEFORE any cell is
#    modified, so cycle rotations don't clobber data we still need to read.
$snapshots = @{}
foreach ($r in $allRows) {
    $rng = $ws.Range("$firstCol$r`:$lastCol$r")
    $snapshots[$r] = $rng.Value2
}

# 2) Force text format on the risky columns for every affected row so the
#    upcoming write-back doesn't get auto-converted by Excel.
foreach ($r in $allRows) {
    foreach ($col in $protectedCols) {
        $ws.Range("$col$r").NumberFormat = "@"
    }
}

# 3) Rotate: write each row's new content from the snapshot of the next row
#    in its cycle.
foreach ($cycle in $cycles) {
    $n = $cycle.Length
    for ($i = 0; $i -lt $n; $i++) {
        $dest = $cycle[$i]
        $src  = $cycle[($i + 1) % $n]
        $ws.Range("$firstCol$dest`:$lastCol$dest").Value2 = $snapshots[$src]
    }
}

# 4) Restore the default General number format on the protected columns
#    (every cell in this sheet normally uses the default/General style).
foreach ($r in $allRows) {
    foreach ($col in $protectedCols) {
        $ws.Range("$col$r").NumberFormat = "General"
    }
}
